$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Header row (row 1) ---
$ws.Range("A1").Value = "Variant (VOC,VOI,VUM)"
$ws.Range("B1").Value = "Other names by which this variant may be known"
$ws.Range("C1").Value = "Lineage/sub-lineages"
$ws.Range("D1").Value = "Number of sequences"
$ws.Range("E1").Value = "Sequences submitted in the last 30 days"
$ws.Range("F1").Value = "date of first sequence"
$ws.Range("G1").Value = "date of last sequence"
$ws.Range("H1").Value = "No of days since last sampling"

# Apply the same bold/centered header style used for A1:C1 to the new G1/H1 header cells
$ws.Range("G1").Font.Bold = $true
$ws.Range("G1").HorizontalAlignment = -4108
$ws.Range("H1").Font.Bold = $true
$ws.Range("H1").HorizontalAlignment = -4108

# Columns F and G (body rows only) hold dates stored as plain text (not Excel
# date serials), matching the source data; force text formatting before
# writing them so the date-like strings aren't auto-converted to numbers.
$ws.Range("F2:G10").NumberFormat = "@"

# --- Clear out the old body (rows 2-10) before rewriting with the new data set/order ---
$ws.Range("A2:H10").ClearContents()

# --- Row 2: A.23.1 ---
$ws.Range("A2").Value = "A.23.1"
$ws.Range("C2").Value = "A.23.1"
$ws.Range("D2").Value = 378
$ws.Range("E2").Value = 0
$ws.Range("F2").Value = "2020-10-21"
$ws.Range("G2").Value = "2021-08-22"
$ws.Range("H2").Value = 39

# --- Row 3: Alpha ---
$ws.Range("A3").Value = "Alpha"
$ws.Range("B3").Value = "VOC-202012/01"
$ws.Range("C3").Value = "B.1.1.7"
$ws.Range("D3").Value = 2498
$ws.Range("E3").Value = 0
$ws.Range("F3").Value = "2020-08-02"
$ws.Range("G3").Value = "2021-08-27"
$ws.Range("H3").Value = 34

# --- Row 4: B.1.1.318 ---
$ws.Range("A4").Value = "B.1.1.318"
$ws.Range("B4").Value = "VUM-2021-06-04"
$ws.Range("C4").Value = "B.1.1.318"
$ws.Range("D4").Value = 465
$ws.Range("E4").Value = 0
$ws.Range("F4").Value = "2021-01-07"
$ws.Range("G4").Value = "2021-08-06"
$ws.Range("H4").Value = 55

# --- Row 5: Beta ---
$ws.Range("A5").Value = "Beta"
$ws.Range("B5").Value = "VOC-202012/02"
$ws.Range("C5").Value = "B.1.351, B.1.351.1, B.1.351.2"
$ws.Range("D5").Value = 10022
$ws.Range("E5").Value = 0
$ws.Range("F5").Value = "2020-05-27"
$ws.Range("G5").Value = "2021-08-27"
$ws.Range("H5").Value = 34

# --- Row 6: C.1 ---
$ws.Range("A6").Value = "C.1"
$ws.Range("C6").Value = "C.1"
$ws.Range("D6").Value = 384
$ws.Range("E6").Value = 0
$ws.Range("F6").Value = "2020-04-12"
$ws.Range("G6").Value = "2021-06-28"
$ws.Range("H6").Value = 94

# --- Row 7: C.1.2 ---
$ws.Range("A7").Value = "C.1.2"
$ws.Range("B7").Value = "VUM-2021-09-01"
$ws.Range("C7").Value = "C.1.2"
$ws.Range("D7").Value = 168
$ws.Range("E7").Value = 7
$ws.Range("F7").Value = "2021-05-11"
$ws.Range("G7").Value = "2021-09-09"
$ws.Range("H7").Value = 21

# --- Row 8: C.36.3 ---
$ws.Range("A8").Value = "C.36.3"
$ws.Range("B8").Value = "VUM-2021-06-16"
$ws.Range("C8").Value = "C.36.3"
$ws.Range("D8").Value = 104
$ws.Range("E8").Value = 0
$ws.Range("F8").Value = "2020-04-26"
$ws.Range("G8").Value = "2021-06-16"
$ws.Range("H8").Value = 106

# --- Row 9: Delta ---
$ws.Range("A9").Value = "Delta"
$ws.Range("B9").Value = "VOC-21APR-02"
$ws.Range("C9").Value = "B.1.617.2, AY.1, AY.3, AY.4, AY.5, AY.6, AY.7.1, AY.10, AY.11, AY.12, AY.13, AY.14, AY.15, AY.16, AY.17, AY.18, AY.19, AY.20, AY.21, AY.23, AY.24, AY.25"
$ws.Range("D9").Value = 12401
$ws.Range("E9").Value = 457
$ws.Range("F9").Value = "2020-09-16"
$ws.Range("G9").Value = "2021-09-20"
$ws.Range("H9").Value = 10

# --- Row 10: Eta ---
$ws.Range("A10").Value = "Eta"
$ws.Range("B10").Value = "VUM-202102/03"
$ws.Range("C10").Value = "B.1.525"
$ws.Range("D10").Value = 960
$ws.Range("E10").Value = 2
$ws.Range("F10").Value = "2020-03-28"
$ws.Range("G10").Value = "2021-09-16"
$ws.Range("H10").Value = 14
